$d = $word.ActiveDocument

# Map of old TOC bookmark names to new TOC bookmark names (renumbered
# after other headings were added earlier in the document's heading
# sequence, per the commit "aggiunta tavola dei contenuti").
$renames = @{
    "_Toc5728339" = "_Toc5728486"
    "_Toc5728340" = "_Toc5728487"
    "_Toc5728341" = "_Toc5728488"
    "_Toc5728342" = "_Toc5728489"
    "_Toc5728343" = "_Toc5728490"
    "_Toc5728344" = "_Toc5728491"
    "_Toc5728345" = "_Toc5728492"
    "_Toc5728346" = "_Toc5728493"
    "_Toc5728347" = "_Toc5728494"
    "_Toc5728348" = "_Toc5728495"
    "_Toc5728349" = "_Toc5728496"
    "_Toc5728350" = "_Toc5728497"
    "_Toc5728351" = "_Toc5728498"
    "_Toc5728352" = "_Toc5728499"
    "_Toc5728353" = "_Toc5728500"
    "_Toc5728354" = "_Toc5728501"
    "_Toc5728355" = "_Toc5728502"
    "_Toc5728356" = "_Toc5728503"
    "_Toc5728357" = "_Toc5728504"
    "_Toc5728358" = "_Toc5728505"
}

foreach ($oldName in $renames.Keys) {
    $newName = $renames[$oldName]
    $bm = $d.Bookmarks.Item($oldName)
    $rng = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($newName, $rng)
}
